# Hexlet progress update of 25/05
# Mark the Transfer-Encoding / query-string / redirects / basic-auth / cookies
# topics on Лист1 as "studied" (value 0, same highlighted style as the
# already-studied rows above them), and fill in the corresponding detail
# notes on Лист2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Лист1")
$ws2 = $wb.Worksheets.Item("Лист2")

# --- Лист1: D13:D17 go from "todo" (1, no fill) to "done" (0, highlighted
# fill) -- copy the format+value from the already-"done" cell D12 onto each
# of them individually (copying onto the whole D13:D17 range at once only
# carries the format, not the value, to every destination cell).
foreach ($r in 13..17) {
    $ws1.Range("D12").Copy($ws1.Cells.Item($r, 4))
}

# --- Лист2: new detail rows 46-68.
$ws2.Range("B46").Value = "Transfer-Encoding"
$ws2.Range("C46").Value = "Transfer-Encoding: chunked"
$ws2.Range("B47").Value = "Передача данных query string"
$ws2.Range("C47").Value = "Строка запроса имеет ограниченную длину"
$ws2.Range("C48").Value = "В POST тоже можно передавать строку запроса, причем вместе с body"
$ws2.Range("C49").Value = "GET идемпотентен, система не меняется"
$ws2.Range("C50").Value = "POST предназначен для изменения данных. Второй ПОСТ может привести к другому ответу и другим модификациям."
$ws2.Range("C51").Value = "ПОСТы никогда не кэшируются"
$ws2.Range("C52").Value = "Формы бывают: - на создание и - на выборку."
$ws2.Range("C53").Value = "на создание отправляется ПОСТ, на выборку - ГЕТ"
$ws2.Range("B54").Value = "Перенаправления"
$ws2.Range("C54").Value = "301 Permanently используется например для перевода с http на https"
$ws2.Range("C55").Value = "Поисковики не выкидывают из индекса такие страницы"
$ws2.Range("C56").Value = "Куды идти указано в заголовке location"
$ws2.Range("B57").Value = "Базовая аутентификация"
$ws2.Range("C57").Value = "Authorization: Basic <base64 encoded login:password>"
$ws2.Range("C58").Value = "Команду printf 'Aladdin:open sesame' | base64 надо, блядь, с кавычками выполнять"
$ws2.Range("B59").Value = "Cookies"
$ws2.Range("C59").Value = "HTTP - stateless protocol. Где тогда чо хранить?"
$ws2.Range("C60").Value = "В curl тоже есть флаги типа --head"
$ws2.Range("C61").Value = "Каждая кука посылается отдельно через свой Set-cookie:"
$ws2.Range("C62").Value = "Кука - это хрень ключ=значение; доп_параметры"
$ws2.Range("C63").Value = "Куки бывают сессионные и персистентные"
$ws2.Range("C64").Value = "domain, path, max-age, expires"
$ws2.Range("C65").Value = "Удаление куки: послать max-age равный НУЛЮ"
$ws2.Range("C66").Value = "HttpOnly - кука не будет отправлена жабаскриптом или аяксом"
$ws2.Range("C67").Value = "Отправка кук: Cookie: key=value; key1=value1"
$ws2.Range("C68").Value = "Параметры кук не отправляются"

# --- Лист2: column B grew a bit wider to fit the new labels.
$ws2.Columns("B").ColumnWidth = 26

# --- Лист2: selection/scroll moved down to the new content. Select on
# Лист2 first (updates its own sheetView selection/scroll), then reactivate
# Лист1 so the workbook's active tab stays on Лист1, matching the diff
# (tabSelected="1" stays on Лист1's sheetView, not Лист2's).
$ws2.Range("B62").Select()
$excel.ActiveWindow.ScrollRow = 41

# --- Лист1: move the selection/scroll down to the newly studied rows, and
# make sure Лист1 ends up as the active tab again.
$ws1.Activate()
$ws1.Range("D19").Select()
$excel.ActiveWindow.ScrollRow = 13
